$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 with text values "123" in both A4 and B4 (matching the
# existing text-typed cells such as B3, not numeric 123).
$a4 = $ws.Range("A4")
$a4.NumberFormat = "@"
$a4.Value = "123"
$a4.ClearFormats()

$b4 = $ws.Range("B4")
$b4.NumberFormat = "@"
$b4.Value = "123"
$b4.ClearFormats()
